$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing "Regex" row (row 4)
# and the following blank row (row 5) down by one, opening up rows 3 and 4
# for the two new config entries.
$ws.Rows.Item(4).Insert()

# New row 3: boolBreakpoint1 config entry
$ws.Range("A3").Value = "boolBreakpoint1"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "breakpoint parameter"

# New row 4: strEmailAccount config entry
$ws.Range("A4").Value = "strEmailAccount"
$ws.Range("B4").Value = "azim.karim@defra.gov.uk"
$ws.Range("C4").Value = "outlook email account"

# Resize the config table to include the two new rows (now A1:C6)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C6"))

# Match the saved selection/active cell
$ws.Range("C4").Select()
